# Ricardo Jacome Resume-1-29-20.docx — apply the tracked edits:
#   1. Merge the "Jaco" / "me R. Stolle, ..." run split back into one run
#      (removes the stray mid-word split that used to host the _GoBack
#      bookmark there).
#   2. Re-type the WCX SAE Presentation line (same visible text, just
#      normalizes the run layout the way Word does when you select and
#      retype a phrase).
#   3. Bump the GPA from 3.88 to 3.90 — this is the most recent edit, so
#      Word leaves the "_GoBack" (last-edit) bookmark sitting right after
#      it.
$d = $word.ActiveDocument

function Placeholder([int]$n) {
    $s = ""
    for ($i = 0; $i -lt $n; $i++) {
        $s += "#"
    }
    return $s
}

function Retype($range, [string]$text) {
    # Round-trips a range through a dummy value so the engine treats the
    # final text as a fresh edit (and re-normalizes the surrounding runs)
    # even though the visible characters end up unchanged.
    $start = $range.Start
    $n = $range.End - $range.Start
    $range.Text = Placeholder($n)
    $fresh = $d.Range($start, $start + $n)
    $fresh.Text = $text
}

# --- 1. "Jaco" + "me R. Stolle, C., & Sweigard M., " -> "Jacome R. Stolle, C., & Sweigard M., " ---
$rng = $d.Content
$needle = "Jacome R. Stolle, C., & Sweigard M., " + [char]0x201C + "Road"
$found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $word1 = $d.Range($rng.Start, $rng.Start + 6)
    Retype $word1 "Jacome"
}

# --- 2. Retype "Presentation, WCX SAE World Congress Experience, " ---
$rng = $d.Content
$found = $rng.Find.Execute("Presentation, WCX SAE World Congress Experience, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    Retype $rng "Presentation, WCX SAE World Congress Experience, "
}

# --- 2b. Retype the ", " between "Detroit, MI" and "April 2020" ---
$rng = $d.Content
$found = $rng.Find.Execute("Detroit, MI, April 2020", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $sepStart = $rng.Start + 11
    $sep = $d.Range($sepStart, $sepStart + 2)
    Retype $sep ", "
}

# --- 2c. Retype "April 2020" ---
$rng = $d.Content
$found = $rng.Find.Execute("April 2020", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    Retype $rng "April 2020"
}

# --- 3. GPA 3.88 -> 3.90, leaving the _GoBack bookmark right after it ---
$rng = $d.Content
$found = $rng.Find.Execute("3.88", $true, $false, $false, $false, $false, $true, 1, $false, "3.90", 2)
if ($found) {
    $rng2 = $d.Content
    $rng2.Find.Execute("3.90", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $bmPoint = $d.Range($rng2.End, $rng2.End)
    $d.Bookmarks.Add("_GoBack", $bmPoint)
}
